# Update from CVSS3.1 to CVSS4.0
# - Re-score three findings (severity text changes)
# - Normalize the finding-ID format for the Open Redirect finding
# - Re-style the Open Redirect severity cell (bold text, darker fill) now
#   that it moved from Low to Medium severity

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Findings")

# Finding #PEN20250001 - XXE in Test Shop: Critical 9.1 -> 9.3
$ws.Range("C6").Value = "Critical (9.3)"

# Finding #PEN20250002 - XSS in Test Shop: High 7.1 -> 7.0
$ws.Range("C7").Value = "High (7.0)"

# Finding previously "#PEN-TEST-0003" - Open Redirect in Test Shop:
# renumber to match the other findings, and rescore Low 3.1 -> Medium 5.6
$ws.Range("A8").Value = "#PEN20250003"
$ws.Range("C8").Value = "Medium (5.6)"

# Highlight the re-scored severity cell: bold font, darker amber/orange fill
$ws.Range("C8").Font.Bold = $true
$ws.Range("C8").Interior.Color = 2468089  # 0x25A8F9 => RGB(0xF9, 0xA8, 0x25) = FFF9A825
